$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: replace the first matchup's data (new teams / new odds) ---
$ws.Range("A2").Value = "OKC"
$ws.Range("B2").Value = "DAL"
$ws.Range("C2").Value = 150
$ws.Range("D2").Value = -180
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 207.5

# --- Rows 5 & 6: fully removed (no games left for these slots) ---
$ws.Rows("5:6").Clear() | Out-Null

# --- Rows 3, 4, 7, 8: team/spread/total data cleared, only the
#     moneyline (C/D) columns remain as blank, styled placeholders ---
$ws.Range("A3:B4").Clear() | Out-Null
$ws.Range("E3:F4").Clear() | Out-Null
$ws.Range("C3:D4").ClearContents() | Out-Null

$ws.Range("A7:B8").Clear() | Out-Null
$ws.Range("E7:F8").Clear() | Out-Null
$ws.Range("C7:D8").ClearContents() | Out-Null

# --- Selection moves from F12 to E12 ---
$ws.Range("E12").Select() | Out-Null
